# Working copy - debug mode new logo
#
# The third "test" sample (row 4) used the shared string "TEST3" in
# column G. It is renamed to "TEST4" (a brand new, previously unused
# string). Because "TEST3" was only referenced by this one cell, giving
# the cell a new value makes "TEST3" unused; Excel drops it from the
# shared-strings table and appends the newly introduced "TEST4" string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "TEST4"

# The sheet's active selection also moved from E11 to G20.
[void]$ws.Range("G20").Select()
